# Updates cryptos list: refresh Price (D) and Volume(1h) (E) columns
# for the rows whose figures changed in this data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.548.02"
$ws.Range("E2").Value = "  -3.93%  "
$ws.Range("D3").Value = "2.540.42"
$ws.Range("E3").Value = "  -3.70%  "
$ws.Range("D5").Value = "506.23"
$ws.Range("E5").Value = "  -4.37%  "
$ws.Range("D6").Value = "143.67"
$ws.Range("E6").Value = "  -7.42%  "
$ws.Range("D7").Value = "0.998"
$ws.Range("D8").Value = "0.565"
$ws.Range("E8").Value = "  -4.15%  "
$ws.Range("D9").Value = "2.545.31"
$ws.Range("E9").Value = "  -3.92%  "
$ws.Range("E10").Value = "  -8.54%  "
$ws.Range("E11").Value = "  -6.19%  "
$ws.Range("D12").Value = "0.332"
$ws.Range("E12").Value = "  -5.51%  "
$ws.Range("E13").Value = "  -0.53%  "
$ws.Range("D14").Value = "2.982.74"
$ws.Range("E14").Value = "  -3.73%  "
$ws.Range("D15").Value = "58.518.60"
$ws.Range("E15").Value = "  -3.96%  "
$ws.Range("E16").Value = "  -5.56%  "
$ws.Range("E17").Value = "  -5.92%  "
$ws.Range("D18").Value = "2.540.23"
$ws.Range("E18").Value = "  -3.97%  "
$ws.Range("E19").Value = "  -4.85%  "
$ws.Range("D20").Value = "339.61"
$ws.Range("E20").Value = "  -3.96%  "
$ws.Range("D21").Value = "10.09"
$ws.Range("E21").Value = "  -5.16%  "
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("E23").Value = "  -4.41%  "
$ws.Range("D24").Value = "60.73"
$ws.Range("E24").Value = "  -1.36%  "
$ws.Range("D25").Value = "0.411"
$ws.Range("E25").Value = "  -4.46%  "
$ws.Range("E26").Value = "  +2.23%  "
$ws.Range("E27").Value = "  -4.96%  "
$ws.Range("D28").Value = "2.652.01"
$ws.Range("E28").Value = "  -3.69%  "
$ws.Range("E29").Value = "  -8.89%  "
$ws.Range("E30").Value = "  -5.99%  "
$ws.Range("E31").Value = "  +0.00%  "
$ws.Range("D32").Value = "149.78"
$ws.Range("E32").Value = "  -0.33%  "
$ws.Range("E33").Value = "  -5.25%  "
$ws.Range("E34").Value = "  -4.88%  "
$ws.Range("D35").Value = "1.53"
$ws.Range("E35").Value = "  -5.54%  "
$ws.Range("D36").Value = "0.915"
$ws.Range("E36").Value = "  +2.38%  "
$ws.Range("D37").Value = "'3.90"
$ws.Range("E37").Value = "  -6.01%  "
$ws.Range("E38").Value = "  -7.28%  "
$ws.Range("D39").Value = "'36.10"
$ws.Range("E39").Value = "  -1.28%  "
$ws.Range("D40").Value = "0.823"
$ws.Range("E40").Value = "  -11.14%  "
$ws.Range("E41").Value = "  -7.04%  "
$ws.Range("D42").Value = "283.02"
$ws.Range("E42").Value = "  -8.32%  "
$ws.Range("E43").Value = "  -7.74%  "
$ws.Range("D44").Value = "0.0998"
$ws.Range("E44").Value = "  -2.25%  "
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("E46").Value = "  -6.41%  "
$ws.Range("E47").Value = "  -5.67%  "
$ws.Range("D48").Value = "18.65"
$ws.Range("E48").Value = "  -5.74%  "
$ws.Range("D49").Value = "10.29"
$ws.Range("E49").Value = "  -0.51%  "
$ws.Range("D50").Value = "0.0226"
$ws.Range("E50").Value = "  -5.13%  "
$ws.Range("E51").Value = "  -9.03%  "
